$d = $word.ActiveDocument

# Locate the paragraph that contains the first chart (immediately after the
# last bullet list item about the order of the series / the lone empty
# paragraph that follows it). The new "Conclusion" section must be inserted
# right before that paragraph.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.InlineShapes.Count -gt 0) {
        $target = $p
        break
    }
}

$targetRange = $target.Range
$null = $targetRange.InsertParagraphBefore()

# Re-fetch the freshly created (still empty) paragraph immediately before the
# chart paragraph and replace its contents with the full Conclusion section
# (heading + body paragraph + two blank paragraphs) via a raw WordprocessingML
# fragment so formatting comes out exactly as authored (no inherited rPr).
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.InlineShapes.Count -gt 0) {
        $insertionPara = $paras.Item($i - 1)
        break
    }
}

$newRange = $insertionPara.Range

$xml = @"
<pkg:xmlData xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:pStyle w:val="Heading1"/>
  </w:pPr>
  <w:r>
    <w:t>Conclusion:</w:t>
  </w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:spacing w:after="0" w:line="276" w:lineRule="auto"/>
  </w:pPr>
  <w:r>
    <w:t>PSM is an efficient and highly adaptive way to simulate ballistic missile trajectory.  By</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>plugging in initial conditions of position, velocity, and acceleration, it can accurately</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>simulate the trajectory of a missile in various situations.  </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">As can be seen in the graphs and data points, the model works as stated.  </w:t>
  </w:r>
  <w:r>
    <w:t>T</w:t>
  </w:r>
  <w:r>
    <w:t>he advantage of being able to take large</w:t>
  </w:r>
  <w:r>
    <w:t>r</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> time steps with PSM</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> compared to RK4 </w:t>
  </w:r>
  <w:r>
    <w:t>mean</w:t>
  </w:r>
  <w:r>
    <w:t>s</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> that any data stored on computers used for solving the simulations are much less, and can be computed faster, while also taking up less space.</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>Overall, PSM can be used in a variety of different applications, as a fast and powerful algorithm.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>
</pkg:xmlData>
"@

$null = $newRange.InsertXML($xml)
